# Weekly "Betarraga" price-history sheet: insert one more historical
# row at the bottom of the series (row 211 -> 212) by shifting the
# Fecha/Volumen/Precio (D, J, K, L, M, P) values of rows 146..211 down
# by one row, and re-seeding row 146 onward so that the data that used
# to sit in the last row (211) now lands in the brand-new row 212 at
# the end of the table. All other columns (A, B, C, E-I, N, O, Q, R)
# are constant for this whole "Betarraga" block, so the new row just
# copies them straight from the row above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 145
$endRow = 211

# --- 1. Snapshot the "before" values for the columns that move ------
$oldD = @{}
$oldJ = @{}
$oldK = @{}
$oldL = @{}
$oldM = @{}
$oldP = @{}

for ($r = $startRow; $r -le $endRow; $r++) {
    $oldD[$r] = $ws.Cells.Item($r, 4).Value2
    $oldJ[$r] = $ws.Cells.Item($r, 10).Value2
    $oldK[$r] = $ws.Cells.Item($r, 11).Value2
    $oldL[$r] = $ws.Cells.Item($r, 12).Value2
    $oldM[$r] = $ws.Cells.Item($r, 13).Value2
    $oldP[$r] = $ws.Cells.Item($r, 16).Value2
}

# --- 2. Shift rows 146..211 down by one source row -------------------
# row r (146..211) takes the old values that used to live in row r-1
for ($r = $startRow + 1; $r -le $endRow; $r++) {
    $src = $r - 1
    $ws.Cells.Item($r, 4).Value  = $oldD[$src]
    $ws.Cells.Item($r, 10).Value = $oldJ[$src]
    $ws.Cells.Item($r, 11).Value = $oldK[$src]
    $ws.Cells.Item($r, 12).Value = $oldL[$src]
    $ws.Cells.Item($r, 13).Value = $oldM[$src]
    $ws.Cells.Item($r, 16).Value = $oldP[$src]
}

# --- 3. Append a brand-new last row (212) with the data that used to
#        be in row 211, duplicating the constant columns from the row
#        right above it and restoring the date's number format (new
#        cells otherwise start out with the default/General style).
$newRow = $endRow + 1

$ws.Cells.Item($newRow, 1).Value  = $ws.Cells.Item($endRow, 1).Value2
$ws.Cells.Item($newRow, 2).Value  = $ws.Cells.Item($endRow, 2).Value2
$ws.Cells.Item($newRow, 3).Value  = $ws.Cells.Item($endRow, 3).Value2

$ws.Cells.Item($newRow, 4).NumberFormat = $ws.Cells.Item($endRow, 4).NumberFormat
$ws.Cells.Item($newRow, 4).Value  = $oldD[$endRow]

$ws.Cells.Item($newRow, 5).Value  = $ws.Cells.Item($endRow, 5).Value2
$ws.Cells.Item($newRow, 6).Value  = $ws.Cells.Item($endRow, 6).Value2
$ws.Cells.Item($newRow, 7).Value  = $ws.Cells.Item($endRow, 7).Value2
$ws.Cells.Item($newRow, 8).Value  = $ws.Cells.Item($endRow, 8).Value2
$ws.Cells.Item($newRow, 9).Value  = $ws.Cells.Item($endRow, 9).Value2

$ws.Cells.Item($newRow, 10).Value = $oldJ[$endRow]
$ws.Cells.Item($newRow, 11).Value = $oldK[$endRow]
$ws.Cells.Item($newRow, 12).Value = $oldL[$endRow]
$ws.Cells.Item($newRow, 13).Value = $oldM[$endRow]

$ws.Cells.Item($newRow, 14).Value = $ws.Cells.Item($endRow, 14).Value2
$ws.Cells.Item($newRow, 15).Value = $ws.Cells.Item($endRow, 15).Value2

$ws.Cells.Item($newRow, 16).Value = $oldP[$endRow]
$ws.Cells.Item($newRow, 17).Value = $ws.Cells.Item($endRow, 17).Value2
$ws.Cells.Item($newRow, 18).Value = $ws.Cells.Item($endRow, 18).Value2
